$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 4 (P03 / 500 / bleeding_broad),
# shifting it down to row 5, then populate the new row 4 with the
# additional bleeding event for P01.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "P01"
$ws.Range("B4").Value = 120
$ws.Range("C4").Value = "bleeding_broad"

$ws.Range("O6").Select()
